# "changed to the report layout"
#
# This script applies the following changes to the activities template:
#  1. Shrink the empty right-aligned paragraph right after "อ.สต.3"
#     (sz 32->12, szCs 32->16) and give it explicit single line spacing.
#  2. Add explicit single line spacing (w:line="240" w:lineRule="auto")
#     to the title paragraph and the "......." signature-line paragraph.
#  3. Collapse the two duplicate empty ListParagraph paragraphs right
#     before "ข้าพเจ้าได้ปฏิบัติงาน..." into a single smaller one
#     (sz 32->16, szCs 32->21).
#  4. Insert three more empty ListParagraph paragraphs right after the
#     {@activitiesRawXml} merge-field paragraph (matching the two that
#     were already there).
#  5. Narrow the page's left/right margins (576 -> 360 twips).
#  6. Reduce the "Normal" style's paragraph spacing-after (200 -> 14 twips).
#  7. Register a new "ListLabel 7" character style (mirrors "ListLabel 6").

$d = $word.ActiveDocument

# --- 1 & 2: line spacing / font-size tweaks on the header paragraphs ---

# Empty paragraph right after "อ.สต.3" (right aligned) - shrink font and
# give it explicit single spacing.
$pNumPage = $d.Paragraphs(4)
$pNumPage.Range.Font.Size = 6
$pNumPage.Range.Font.SizeBi = 8
$pNumPage.Format.LineSpacingRule = 0

# Title paragraph ("แบบรายงานการปฏิบัติงาน...") - just add spacing.
$d.Paragraphs(5).Format.LineSpacingRule = 0

# "......." signature-line paragraph - just add spacing.
$d.Paragraphs(6).Format.LineSpacingRule = 0

# --- 3: shrink the first duplicate empty paragraph, drop the second ---

$pDup1 = $d.Paragraphs(15)
$pDup1.Range.Font.Size = 8
$pDup1.Range.Font.SizeBi = 10.5

$d.Paragraphs(16).Range.Delete()

# --- 4: insert three more empty ListParagraph paragraphs after the
#        {@activitiesRawXml} paragraph (now shifted to index 17) ---

$pActivities = $d.Paragraphs(17)
$pActivities.Range.InsertParagraphAfter()
$pActivities.Range.InsertParagraphAfter()
$pActivities.Range.InsertParagraphAfter()

for ($i = 18; $i -le 20; $i++) {
    $newPara = $d.Paragraphs($i)
    $newPara.Range.Font.SizeBi = 16
}

# --- 5: narrow the page margins ---

$d.PageSetup.LeftMargin = 18
$d.PageSetup.RightMargin = 18

# --- 6: reduce "Normal" style spacing-after ---

$d.Styles("Normal").ParagraphFormat.SpaceAfter = 0.7

# --- 7: register the new "ListLabel 7" character style ---

$listLabel7 = $d.Styles.Add("ListLabel7", 2)
$listLabel7.NameLocal = "ListLabel 7"
$listLabel7.Font.Size = 2
$listLabel7.Font.SizeBi = 2.5
